{"js": "// Replace the three-digit-by-one-digit multiplication expressions in the\n// table cells with a new set of problems (same \"A\u00d7B=C\" text format).\n// Each old string occurs exactly once in the document, so a plain\n// search-and-replace per pair is unambiguous.\nconst replacements = [\n  [\"752\u00d77=5264\", \"489\u00d74=1956\"],\n  [\"988\u00d78=7904\", \"434\u00d77=3038\"],\n  [\"444\u00d78=3552\", \"212\u00d77=1484\"],\n  [\"892\u00d75=4460\", \"780\u00d76=4680\"],\n  [\"998\u00d74=3992\", \"881\u00d77=6167\"],\n  [\"389\u00d76=2334\", \"650\u00d77=4550\"],\n  [\"278\u00d77=1946\", \"376\u00d77=2632\"],\n  [\"943\u00d78=7544\", \"925\u00d73=2775\"],\n  [\"494\u00d72=988\", \"755\u00d73=2265\"],\n  [\"432\u00d79=3888\", \"806\u00d76=4836\"],\n  [\"260\u00d72=520\", \"974\u00d76=5844\"],\n  [\"788\u00d72=1576\", \"626\u00d76=3756\"],\n  [\"466\u00d76=2796\", \"261\u00d76=1566\"],\n  [\"643\u00d78=5144\", \"457\u00d77=3199\"],\n  [\"369\u00d74=1476\", \"541\u00d76=3246\"],\n  [\"866\u00d76=5196\", \"975\u00d78=7800\"],\n  [\"895\u00d78=7160\", \"961\u00d73=2883\"],\n  [\"983\u00d75=4915\", \"796\u00d72=1592\"],\n  [\"335\u00d76=2010\", \"914\u00d74=3656\"],\n  [\"786\u00d74=3144\", \"467\u00d77=3269\"],\n  [\"878\u00d79=7902\", \"136\u00d73=408\"],\n  [\"582\u00d78=4656\", \"693\u00d79=6237\"],\n  [\"480\u00d72=960\", \"648\u00d75=3240\"],\n  [\"246\u00d78=1968\", \"406\u00d76=2436\"],\n  [\"274\u00d79=2466\", \"829\u00d73=2487\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`No match found for \"${oldText}\"`);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the three-digit-by-one-digit multiplication expressions in the\n# table cells with a new set of problems (same \"A\u00d7B=C\" text format).\n# Each old string occurs exactly once in the document, so Find/Replace on\n# the whole document body is unambiguous for each pair.\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = '752\u00d77=5264'; New = '489\u00d74=1956' },\n    @{ Old = '988\u00d78=7904'; New = '434\u00d77=3038' },\n    @{ Old = '444\u00d78=3552'; New = '212\u00d77=1484' },\n    @{ Old = '892\u00d75=4460'; New = '780\u00d76=4680' },\n    @{ Old = '998\u00d74=3992'; New = '881\u00d77=6167' },\n    @{ Old = '389\u00d76=2334'; New = '650\u00d77=4550' },\n    @{ Old = '278\u00d77=1946'; New = '376\u00d77=2632' },\n    @{ Old = '943\u00d78=7544'; New = '925\u00d73=2775' },\n    @{ Old = '494\u00d72=988';  New = '755\u00d73=2265' },\n    @{ Old = '432\u00d79=3888'; New = '806\u00d76=4836' },\n    @{ Old = '260\u00d72=520';  New = '974\u00d76=5844' },\n    @{ Old = '788\u00d72=1576'; New = '626\u00d76=3756' },\n    @{ Old = '466\u00d76=2796'; New = '261\u00d76=1566' },\n    @{ Old = '643\u00d78=5144'; New = '457\u00d77=3199' },\n    @{ Old = '369\u00d74=1476'; New = '541\u00d76=3246' },\n    @{ Old = '866\u00d76=5196'; New = '975\u00d78=7800' },\n    @{ Old = '895\u00d78=7160'; New = '961\u00d73=2883' },\n    @{ Old = '983\u00d75=4915'; New = '796\u00d72=1592' },\n    @{ Old = '335\u00d76=2010'; New = '914\u00d74=3656' },\n    @{ Old = '786\u00d74=3144'; New = '467\u00d77=3269' },\n    @{ Old = '878\u00d79=7902'; New = '136\u00d73=408'  },\n    @{ Old = '582\u00d78=4656'; New = '693\u00d79=6237' },\n    @{ Old = '480\u00d72=960';  New = '648\u00d75=3240' },\n    @{ Old = '246\u00d78=1968'; New = '406\u00d76=2436' },\n    @{ Old = '274\u00d79=2466'; New = '829\u00d73=2487' }\n)\n\nforeach ($pair in $replacements) {\n    $rng = $d.Content\n    $found = $rng.Find.Execute(\n        $pair.Old,  # FindText\n        $false,     # MatchCase\n        $false,     # MatchWholeWord\n        $false,     # MatchWildcards\n        $false,     # MatchSoundsLike\n        $false,     # MatchAllWordForms\n        $true,      # Forward\n        1,          # Wrap (wdFindContinue)\n        $false,     # Format\n        $pair.New,  # ReplaceWith\n        2           # Replace (wdReplaceAll)\n    )\n    if (-not $found) {\n        throw \"No match found for '$($pair.Old)'\"\n    }\n}\n"}
